$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new rows at 345-349, shifting existing rows 345-362 down to 350-367
$ws.Range("A345:A349").EntireRow.Insert()

# Row 345
$ws.Cells.Item(345, 1).Value = 11
$ws.Cells.Item(345, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(345, 3).Value = "Bíobío"
$ws.Cells.Item(345, 4).Value = 44516
$ws.Cells.Item(345, 5).Value = 8
$ws.Cells.Item(345, 6).Value = "Fruta"
$ws.Cells.Item(345, 7).Value = 100106
$ws.Cells.Item(345, 8).Value = "Oleaginosos"
$ws.Cells.Item(345, 9).Value = 100106002
$ws.Cells.Item(345, 10).Value = "Palta"
$ws.Cells.Item(345, 11).Value = "Edranol"
$ws.Cells.Item(345, 12).Value = "Primera"
$ws.Cells.Item(345, 13).Value = 50
$ws.Cells.Item(345, 14).Value = 2200
$ws.Cells.Item(345, 15).Value = 2200
$ws.Cells.Item(345, 16).Value = 2200
$ws.Cells.Item(345, 17).Value = "$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(345, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(345, 19).Value = 2200
$ws.Cells.Item(345, 20).Value = 1

# Row 346
$ws.Cells.Item(346, 1).Value = 11
$ws.Cells.Item(346, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(346, 3).Value = "Bíobío"
$ws.Cells.Item(346, 4).Value = 44516
$ws.Cells.Item(346, 5).Value = 8
$ws.Cells.Item(346, 6).Value = "Fruta"
$ws.Cells.Item(346, 7).Value = 100106
$ws.Cells.Item(346, 8).Value = "Oleaginosos"
$ws.Cells.Item(346, 9).Value = 100106002
$ws.Cells.Item(346, 10).Value = "Palta"
$ws.Cells.Item(346, 11).Value = "Edranol"
$ws.Cells.Item(346, 12).Value = "Segunda"
$ws.Cells.Item(346, 13).Value = 50
$ws.Cells.Item(346, 14).Value = 2000
$ws.Cells.Item(346, 15).Value = 2000
$ws.Cells.Item(346, 16).Value = 2000
$ws.Cells.Item(346, 17).Value = "$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(346, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(346, 19).Value = 2000
$ws.Cells.Item(346, 20).Value = 1

# Row 347
$ws.Cells.Item(347, 1).Value = 11
$ws.Cells.Item(347, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(347, 3).Value = "Bíobío"
$ws.Cells.Item(347, 4).Value = 44516
$ws.Cells.Item(347, 5).Value = 8
$ws.Cells.Item(347, 6).Value = "Fruta"
$ws.Cells.Item(347, 7).Value = 100106
$ws.Cells.Item(347, 8).Value = "Oleaginosos"
$ws.Cells.Item(347, 9).Value = 100106002
$ws.Cells.Item(347, 10).Value = "Palta"
$ws.Cells.Item(347, 11).Value = "Hass"
$ws.Cells.Item(347, 12).Value = "Primera"
$ws.Cells.Item(347, 13).Value = 50
$ws.Cells.Item(347, 14).Value = 2800
$ws.Cells.Item(347, 15).Value = 2800
$ws.Cells.Item(347, 16).Value = 2800
$ws.Cells.Item(347, 17).Value = "$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(347, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(347, 19).Value = 2800
$ws.Cells.Item(347, 20).Value = 1

# Row 348
$ws.Cells.Item(348, 1).Value = 11
$ws.Cells.Item(348, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(348, 3).Value = "Bíobío"
$ws.Cells.Item(348, 4).Value = 44516
$ws.Cells.Item(348, 5).Value = 8
$ws.Cells.Item(348, 6).Value = "Fruta"
$ws.Cells.Item(348, 7).Value = 100106
$ws.Cells.Item(348, 8).Value = "Oleaginosos"
$ws.Cells.Item(348, 9).Value = 100106002
$ws.Cells.Item(348, 10).Value = "Palta"
$ws.Cells.Item(348, 11).Value = "Hass"
$ws.Cells.Item(348, 12).Value = "Segunda"
$ws.Cells.Item(348, 13).Value = 50
$ws.Cells.Item(348, 14).Value = 2600
$ws.Cells.Item(348, 15).Value = 2600
$ws.Cells.Item(348, 16).Value = 2600
$ws.Cells.Item(348, 17).Value = "$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(348, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(348, 19).Value = 2600
$ws.Cells.Item(348, 20).Value = 1

# Row 349
$ws.Cells.Item(349, 1).Value = 11
$ws.Cells.Item(349, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(349, 3).Value = "Bíobío"
$ws.Cells.Item(349, 4).Value = 44516
$ws.Cells.Item(349, 5).Value = 8
$ws.Cells.Item(349, 6).Value = "Fruta"
$ws.Cells.Item(349, 7).Value = 100106
$ws.Cells.Item(349, 8).Value = "Oleaginosos"
$ws.Cells.Item(349, 9).Value = 100106002
$ws.Cells.Item(349, 10).Value = "Palta"
$ws.Cells.Item(349, 11).Value = "Hass"
$ws.Cells.Item(349, 12).Value = "Tercera"
$ws.Cells.Item(349, 13).Value = 50
$ws.Cells.Item(349, 14).Value = 2300
$ws.Cells.Item(349, 15).Value = 2300
$ws.Cells.Item(349, 16).Value = 2300
$ws.Cells.Item(349, 17).Value = "$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(349, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(349, 19).Value = 2300
$ws.Cells.Item(349, 20).Value = 1
